$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "TTFF AOI (ms)" row entirely (row 5). Deleting the row
# shifts rows 6-11 up by one, so the remaining metric rows line up with
# their new row numbers (Dwell time (ms) -> 5, Dwell time (%) -> 6,
# Fixation duration (ms) -> 7, First fixation duration (ms) -> 8).
$ws.Rows("5:5").Delete()

# --- Drop the two now-trailing blank rows (formerly rows 10 & 11, now
# rows 9 & 10 after the shift above) so the used range ends at row 8.
$ws.Rows("9:10").Delete()

# --- Clear the header-row formatting (bold font / border / centered
# alignment) so the header cells fall back to the default style.
$ws.Range("A1:P1").Style = "Normal"

# --- The "Unnamed: 0" label in A1 is cleared to an empty string.
$ws.Range("A1").Value = ""

# --- Corrected data values (pre/post/total fixation cleaning fixes).
# Row 3: Revisit count
$ws.Range("B3").Value = 11
$ws.Range("K3").Value = 4

# Row 4: Fixation count
$ws.Range("B4").Value = 32
$ws.Range("K4").Value = 150

# Row 5: Dwell time (ms)
$ws.Range("B5").Value = 8559.360000000001
$ws.Range("K5").Value = 43981.63

# Row 6: Dwell time (%)
$ws.Range("B6").Value = 12.62
$ws.Range("C6").Value = 0.52
$ws.Range("E6").Value = 0.52
$ws.Range("F6").Value = 0.52
$ws.Range("G6").Value = 2.95
$ws.Range("J6").Value = 6.4
$ws.Range("K6").Value = 64.84999999999999
$ws.Range("L6").Value = 0.34
$ws.Range("M6").Value = 3.3

# Row 7: Fixation duration (ms)
$ws.Range("B7").Value = 267.48
$ws.Range("K7").Value = 293.21
